{"js": "// Convert US-style decimal numbers (1,234.56) to European-style (1.234,56)\n// everywhere in the document body, including inside table cells.\n// Dates (DD.MM.YYYY) and other text are left untouched because they do not\n// match the \"digits[,digits]*.digits(2)\" shape.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Matches strings like \"71.12\", \"5,971.44\", \"1,234,567.89\"\nconst usNumberRe = /^\\d{1,3}(,\\d{3})*\\.\\d{2}$/;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (usNumberRe.test(text)) {\n    const segments = text.split(\".\");\n    const decimals = segments.pop();\n    const integerPart = segments.join(\".\").split(\",\").join(\".\");\n    const europeanText = integerPart + \",\" + decimals;\n    if (europeanText !== text) {\n      para.getRange().insertText(europeanText, \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Convert US-style decimal numbers (1,234.56) to European-style (1.234,56)\n# everywhere in the document, including inside table cells.\n# Dates (DD.MM.YYYY) and other text are left untouched because they do not\n# match the \"digits[,digits]*.digits(2)\" shape.\n\n$d = $word.ActiveDocument\n\n# Each paragraph's Range.Text ends with a trailing mark: \"\\r\" for a normal\n# paragraph, or \"\\r\\a\" for the last paragraph inside a table cell. Allow for\n# either (or none) after the number.\n$pattern = '^(\\d{1,3}(,\\d{3})*\\.\\d{2})[\\r\\a]*$'\n\n$count = 0\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $match = [regex]::Match($t, $pattern)\n    if ($match.Success) {\n        $orig = $match.Groups[1].Value\n\n        # Swap thousands separator \",\" and decimal separator \".\" :\n        # split off the decimal part, turn remaining \",\" into \".\" for the\n        # thousands groups, then join with \",\" as the new decimal separator.\n        $parts = $orig.Split('.')\n        $dec = $parts[$parts.Length - 1]\n        $intPart = ($parts[0..($parts.Length - 2)] -join '.').Replace(',', '.')\n        $euro = \"$intPart,$dec\"\n\n        if ($euro -ne $orig) {\n            $r = $p.Range\n            $r.SetRange($r.Start, $r.Start + $orig.Length)\n            $r.Text = $euro\n            $count++\n        }\n    }\n}\n\nWrite-Output \"Replaced $count number(s)\"\n"}
